# Adjusted to Naming Convention: rename EPP* process-set identifiers to P-TH*
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC_CCS")

# Column G (rows 3-8): EPP*CCS* -> P-TH*CCS*
for ($r = 3; $r -le 8; $r++) {
    $ws.Cells.Item($r, 7).Value = "P-TH*CCS*"
}

# Column C (rows 16-20): EPP*<Area>* -> P-TH*<Area>*
$ws.Range("C16").Value = "P-TH*Moneypoint*"
$ws.Range("C17").Value = "P-TH*Cork*"
$ws.Range("C18").Value = "P-TH*Dublin*"
$ws.Range("C19").Value = "P-TH*Offaly*"
$ws.Range("C20").Value = "P-TH*Kilroot*"
